# Update column F (dSF) values on the active worksheet to reflect
# repulled/recalculated data, per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = 2
    4  = -2
    5  = 1
    8  = -2
    9  = -2
    10 = 9
    11 = 1
    12 = -2
    13 = 2
    14 = 2
    15 = 5
    16 = -1
    17 = 5
    18 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
